# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and fix a row-order swap between BitcoinCash and InternetComputer(DFINITY)
# (rows 24/25).
#
# Values that look numeric (e.g. "413.34", "1.00") must be written as TEXT
# (the sheet stores Price/Volume as strings), so we prefix them with a
# leading apostrophe to force text entry, then reset Style to "Normal" so
# the cell keeps the workbook's default (unstyled) formatting rather than
# picking up a "quote prefix" text number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.592.57"
$ws.Range("E2").Value = "  +6.83%  "
$ws.Range("D3").Value = "3.469.42"
$ws.Range("E3").Value = "  +5.07%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'413.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("D6").Value = "'128.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +17.49%  "
$ws.Range("D7").Value = "3.460.95"
$ws.Range("E7").Value = "  +4.94%  "
$ws.Range("D8").Value = "'0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'0.694"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.26%  "
$ws.Range("E11").Value = "  +29.48%  "
$ws.Range("D12").Value = "'42.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.33%  "
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "4.031.97"
$ws.Range("E14").Value = "  +5.32%  "
$ws.Range("D15").Value = "'8.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.04%  "
$ws.Range("D16").Value = "'20.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.79%  "
$ws.Range("D17").Value = "3.480.37"
$ws.Range("E17").Value = "  +5.34%  "
$ws.Range("D18").Value = "62.618.53"
$ws.Range("E18").Value = "  +7.26%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "'10.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").Value = "'0.0000135"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +23.74%  "
$ws.Range("D22").Value = "'3.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").Value = "'82.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.81%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'13.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("B25").Value = "BitcoinCash"
$ws.Range("C25").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D25").Value = "'314.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.47%  "
$ws.Range("D26").Value = "'3.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'30.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.17%  "
$ws.Range("D28").Value = "'8.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.22%  "
$ws.Range("D29").Value = "'7.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.92%  "
$ws.Range("E30").Value = "  +5.17%  "
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("E32").Value = "  +4.24%  "
$ws.Range("D33").Value = "'2.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +23.51%  "
$ws.Range("D34").Value = "'11.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("D35").Value = "'42.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.76%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  -5.50%  "
$ws.Range("D38").Value = "'52.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").Value = "'3.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").Value = "'3.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.05%  "
$ws.Range("D42").Value = "'2.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.03%  "
$ws.Range("D43").Value = "'0.126"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.60%  "
$ws.Range("D44").Value = "'136.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "'17.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("D46").Value = "'0.286"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D47").Value = "'3.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "'2.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("D49").Value = "'22.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("D50").Value = "2.225.85"
$ws.Range("E50").Value = "  +2.82%  "
$ws.Range("D51").Value = "3.825.79"
$ws.Range("E51").Value = "  +5.36%  "
